$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.637.04"
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$ws.Range("D3").Value = "1.798.78"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "226.90"
$r.Style = "Normal"

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.558"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "

# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "32.98"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +3.62%  "

# Row 9
$ws.Range("E9").Value = "  +2.22%  "

# Row 10
$ws.Range("E10").Value = "  +1.26%  "

# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0950"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "

# Row 12
$ws.Range("D12").Value = "2.056.34"
$ws.Range("E12").Value = "  +0.63%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.822.00"
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "11.17"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.640"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +2.72%  "

# Row 16
$ws.Range("D16").Value = "34.583.48"
$ws.Range("E16").Value = "  +1.47%  "

# Row 17
$ws.Range("E17").Value = "  +2.72%  "

# Row 18
$ws.Range("E18").Value = "  +1.10%  "

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "248.40"
$r.Style = "Normal"

# Row 20
$ws.Range("E20").Value = "  +3.36%  "

# Row 21
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "11.33"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +3.93%  "

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "4.18"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("E24").Value = "  +0.33%  "

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "164.89"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +2.13%  "

# Row 26
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "7.27"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "

# Row 27
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "16.56"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.117"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +2.73%  "

# Row 29
$ws.Range("E29").Value = "  -0.23%  "

# Row 30
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "3.99"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +9.55%  "

# Row 31
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.82"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +3.66%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "1.24"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.0523"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "

# Row 34
$ws.Range("E34").Value = "  +1.87%  "

# Row 35
$ws.Range("D35").Value = "1.425.92"
$ws.Range("E35").Value = "  -1.63%  "

# Row 36
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "2.57"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +5.58%  "

# Row 37
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.674"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +2.86%  "

# Row 38
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.0193"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("E39").Value = "  +1.97%  "

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "85.53"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +6.32%  "

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.939"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +1.86%  "

# Row 42
$ws.Range("E42").Value = "  +0.61%  "

# Row 43
$ws.Range("E43").Value = "  +2.30%  "

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "13.53"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.0524"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +2.82%  "

# Row 46
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "6.09"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "

# Row 47
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.08"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
$ws.Range("D48").Value = "1.956.12"
$ws.Range("E48").Value = "  +0.56%  "

# Row 49
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "106.23"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "

# Row 51
$ws.Range("E51").Value = "  -4.37%  "
